# Apply the VisaoFuncional_FluxosdeEventos text corrections described by
# the commit:
#   "docs: correcao da visao de funcao ... correcao do visao de funcao
#    docx, adicao de erro no login ao deixar campo sem preecnher"
#
# All edits are scoped to individual table cells (via Range.Start/End
# arithmetic) rather than whole-document Find/Replace, because
# Find.Execute on this host always searches from the top of the document
# regardless of which Range's .Find was invoked - a whole-document
# replace would incorrectly touch the "Subsistema/modulo:" title
# paragraph (which the diff leaves untouched) as well as every other
# occurrence of the shared phrase across the 11 scenario tables.

$d = $word.ActiveDocument

function Set-CellText {
    # Replace the *entire* contents of a table cell (collapses it to a
    # single paragraph/run with the surrounding formatting).
    param($table, $row, $col, [string]$newText)
    $cell = $table.Cell($row, $col)
    $start = $cell.Range.Start
    $end = $cell.Range.End
    $r = $d.Range($start, $end - 1)
    $r.Text = $newText
}

function Set-CellSuffix {
    # Replace everything in a (single-paragraph) cell *after* a literal
    # prefix, leaving the prefix run untouched.
    param($table, $row, $col, [string]$prefix, [string]$newSuffix)
    $cell = $table.Cell($row, $col)
    $start = $cell.Range.Start
    $end = $cell.Range.End
    $prefixLen = $prefix.Length
    $r = $d.Range($start + $prefixLen, $end - 1)
    $r.Text = $newSuffix
}

function Add-TextAfterParagraph {
    # Append text right before a paragraph's end-of-paragraph mark
    # (used to extend the last run of a paragraph without disturbing
    # the rest of the cell, e.g. cells with multiple paragraphs).
    param($paragraph, [string]$extra)
    $p = $paragraph.Range
    $insertPoint = $d.Range($p.End - 1, $p.End - 1)
    $insertPoint.InsertAfter($extra)
}

# ---------------------------------------------------------------------
# Table 1 ("Modulo Autenticacao" / cenario correto)
# ---------------------------------------------------------------------
$t1 = $d.Tables(1)
Set-CellText   $t1 1 2 "Realizar Autenticação"
Set-CellText   $t1 2 2 "Realizar Autenticação Correto"
Set-CellText   $t1 3 2 "Pessoas"
Set-CellSuffix $t1 5 2 "O " "Usuário autentica-se no sistema"

# ---------------------------------------------------------------------
# Table 2 ("Modulo Autenticacao Senha Incorreta")
# ---------------------------------------------------------------------
$t2 = $d.Tables(2)
Set-CellText $t2 1 2 "Realizar Autenticação"
Set-CellText $t2 2 2 "Realizar Autenticação senha incorreta"
Set-CellText $t2 3 2 "Pessoas"
Set-CellText $t2 5 2 "Mensagem de erro e-mail ou senha inválidos"

# ---------------------------------------------------------------------
# Table 3 ("Modulo Autenticacao E-mail Incorreto")
# ---------------------------------------------------------------------
$t3 = $d.Tables(3)
Set-CellText $t3 1 2 "Realizar Autenticação"
Set-CellText $t3 2 2 "Realizar Autenticação E-mail Incorreto"
Set-CellText $t3 3 2 "Pessoas"
Set-CellText $t3 5 2 "Mensagem de erro e-mail ou senha inválidos"

# ---------------------------------------------------------------------
# Table 4 ("Modulo Autenticacao campos sem preencher")
# ---------------------------------------------------------------------
$t4 = $d.Tables(4)
Set-CellText $t4 1 2 "Realizar Autenticação"
Set-CellText $t4 2 2 "Realizar Autenticação campos sem preencher"
Set-CellText $t4 3 2 "Pessoas"
Set-CellText $t4 5 2 "Mensagem de erro, por favor, preencha todos os campos obrigatórios."

# ---------------------------------------------------------------------
# Table 5 ("Acessar Recuperar Senha" correto)
# ---------------------------------------------------------------------
$t5 = $d.Tables(5)
Set-CellText $t5 3 2 "Pessoas"

# ---------------------------------------------------------------------
# Table 6 ("Acessar Recuperar Senha" e-mail incorreto)
# ---------------------------------------------------------------------
$t6 = $d.Tables(6)
Set-CellText $t6 5 2 "Mensagem de erro do sistema"
$t6cell = $t6.Cell(10, 2)
Add-TextAfterParagraph $t6cell.Range.Paragraphs(1) " que não foi possível enviar o código"

# ---------------------------------------------------------------------
# Table 7 ("Cadastro de Usuario no Sistema" correto)
# ---------------------------------------------------------------------
$t7 = $d.Tables(7)
Set-CellText $t7 3 2 "Pessoas"

# ---------------------------------------------------------------------
# Tables 8-11: only the "Ator(es) envolvido(s)" cell loses ", sistema"
# ---------------------------------------------------------------------
Set-CellText $d.Tables(8)  3 2 "Pessoas"
Set-CellText $d.Tables(9)  3 2 "Pessoas"
Set-CellText $d.Tables(10) 3 2 "Pessoas"
Set-CellText $d.Tables(11) 3 2 "Pessoas"

Write-Output "ok"
